$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.357.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.794.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3806"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3453"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.205"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07526"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.500"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.792.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.087"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001100"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06663"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.539"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.356.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.430"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.585"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.502"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "152.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.997.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "134.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.064"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.142"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08722"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.689"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.461"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6917"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.900"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06382"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2207"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02338"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.276"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6483"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.868"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.128"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "130.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07194"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.93%  "
